$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 68
$ws.Range("H2").Value = "bedrooms"
$ws.Range("I2").Value = "distractor"
$ws.Range("K2").Value = "f"
$ws.Range("L2").Value = "stimuli/img_i7vab.png"
$ws.Range("M2").Value = 86.40000000000001
$ws.Range("N2").Value = 67.8
$ws.Range("O2").Value = 77.09999999999999
$ws.Range("P2").Value = 35
$ws.Range("F3").Value = 69
$ws.Range("H3").Value = "kitchens"
$ws.Range("I3").Value = "target"
$ws.Range("K3").Value = "j"
$ws.Range("L3").Value = "stimuli/img_84s7n.png"
$ws.Range("M3").Value = 11.03125
$ws.Range("N3").Value = 2.90625
$ws.Range("O3").Value = 6.96875
$ws.Range("P3").Value = 32
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 1
$ws.Range("F4").Value = 70
$ws.Range("L4").Value = "stimuli/img_411xa.png"
$ws.Range("M4").Value = 51.03030303030303
$ws.Range("N4").Value = 28.93939393939394
$ws.Range("O4").Value = 39.98484848484848
$ws.Range("P4").Value = 33
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 2
$ws.Range("F5").Value = 71
$ws.Range("H5").Value = "bedrooms"
$ws.Range("L5").Value = "stimuli/img_scrdm.png"
$ws.Range("M5").Value = 78.675
$ws.Range("N5").Value = 57.9
$ws.Range("O5").Value = 68.28749999999999
$ws.Range("F6").Value = 72
$ws.Range("L6").Value = "stimuli/img_uy1n4.png"
$ws.Range("M6").Value = 76.30555555555556
$ws.Range("N6").Value = 55.33333333333334
$ws.Range("O6").Value = 65.81944444444444
$ws.Range("P6").Value = 36
$ws.Range("Q6").Value = 7
$ws.Range("R6").Value = 7
$ws.Range("S6").Value = 7
$ws.Range("F7").Value = 73
$ws.Range("L7").Value = "stimuli/img_xti0z.png"
$ws.Range("M7").Value = 81.40625
$ws.Range("N7").Value = 61.4375
$ws.Range("O7").Value = 71.421875
$ws.Range("P7").Value = 32
$ws.Range("Q7").Value = 8
$ws.Range("R7").Value = 8
$ws.Range("S7").Value = 8
$ws.Range("F8").Value = 74
$ws.Range("L8").Value = "stimuli/img_jz3kd.png"
$ws.Range("M8").Value = 72.79411764705883
$ws.Range("N8").Value = 51.64705882352941
$ws.Range("O8").Value = 62.22058823529412
$ws.Range("P8").Value = 34
$ws.Range("Q8").Value = 6
$ws.Range("R8").Value = 6
$ws.Range("S8").Value = 6
$ws.Range("F9").Value = 75
$ws.Range("F10").Value = 76
$ws.Range("L10").Value = "stimuli/img_7pgd2.png"
$ws.Range("M10").Value = 78.59375
$ws.Range("N10").Value = 57.84375
$ws.Range("O10").Value = 68.21875
$ws.Range("P10").Value = 32
$ws.Range("Q10").Value = 8
$ws.Range("R10").Value = 7
$ws.Range("S10").Value = 7
$ws.Range("F11").Value = 77
$ws.Range("L11").Value = "stimuli/img_j5rpx.png"
$ws.Range("M11").Value = 72.24242424242425
$ws.Range("N11").Value = 50
$ws.Range("O11").Value = 61.12121212121212
$ws.Range("P11").Value = 33
$ws.Range("Q11").Value = 5
$ws.Range("R11").Value = 5
$ws.Range("S11").Value = 5
$ws.Range("F12").Value = 78
$ws.Range("L12").Value = "stimuli/img_6nbgt.png"
$ws.Range("M12").Value = 78.45161290322581
$ws.Range("N12").Value = 57.83870967741935
$ws.Range("O12").Value = 68.14516129032258
$ws.Range("P12").Value = 31
$ws.Range("Q12").Value = 7
$ws.Range("F13").Value = 79
$ws.Range("H13").Value = "bedrooms"
$ws.Range("I13").Value = "distractor"
$ws.Range("K13").Value = "f"
$ws.Range("L13").Value = "stimuli/img_bklr1.png"
$ws.Range("M13").Value = 86.54761904761905
$ws.Range("N13").Value = 67.73809523809524
$ws.Range("O13").Value = 77.14285714285714
$ws.Range("P13").Value = 42
$ws.Range("Q13").Value = 9
$ws.Range("F14").Value = 80
$ws.Range("L14").Value = "stimuli/img_c0me7.png"
$ws.Range("M14").Value = 68.40000000000001
$ws.Range("N14").Value = 45.62857142857143
$ws.Range("O14").Value = 57.01428571428572
$ws.Range("Q14").Value = 4
$ws.Range("R14").Value = 4
$ws.Range("S14").Value = 4
$ws.Range("F15").Value = 81
$ws.Range("L15").Value = "stimuli/img_xesl0.png"
$ws.Range("M15").Value = 69.28571428571429
$ws.Range("N15").Value = 47.35714285714285
$ws.Range("O15").Value = 58.32142857142857
$ws.Range("P15").Value = 28
$ws.Range("Q15").Value = 5
$ws.Range("R15").Value = 5
$ws.Range("S15").Value = 5
$ws.Range("F16").Value = 82
$ws.Range("L16").Value = "stimuli/img_cv9qj.png"
$ws.Range("M16").Value = 60.34375
$ws.Range("N16").Value = 35.34375
$ws.Range("O16").Value = 47.84375
$ws.Range("P16").Value = 32
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = 3
$ws.Range("S16").Value = 3
$ws.Range("F17").Value = 83
$ws.Range("H17").Value = "living_rooms"
$ws.Range("I17").Value = "distractor"
$ws.Range("K17").Value = "f"
$ws.Range("L17").Value = "stimuli/img_24rt2.png"
$ws.Range("M17").Value = 55.26829268292683
$ws.Range("N17").Value = 34.19512195121951
$ws.Range("O17").Value = 44.73170731707317
$ws.Range("P17").Value = 41
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = 3
$ws.Range("S17").Value = 3
$ws.Range("F18").Value = 84
$ws.Range("L18").Value = "stimuli/img_ncr40.png"
$ws.Range("M18").Value = 75.66666666666667
$ws.Range("N18").Value = 54.27272727272727
$ws.Range("O18").Value = 64.96969696969697
$ws.Range("P18").Value = 33
$ws.Range("Q18").Value = 6
$ws.Range("R18").Value = 6
$ws.Range("S18").Value = 6
$ws.Range("F19").Value = 85
$ws.Range("H19").Value = "kitchens"
$ws.Range("I19").Value = "target"
$ws.Range("K19").Value = "j"
$ws.Range("L19").Value = "stimuli/img_s9are.png"
$ws.Range("M19").Value = 90.14285714285714
$ws.Range("N19").Value = 75.22857142857143
$ws.Range("O19").Value = 82.68571428571428
$ws.Range("P19").Value = 35
$ws.Range("Q19").Value = 10
$ws.Range("R19").Value = 10
$ws.Range("S19").Value = 10
$ws.Range("F20").Value = 86
$ws.Range("H20").Value = "kitchens"
$ws.Range("I20").Value = "target"
$ws.Range("K20").Value = "j"
$ws.Range("L20").Value = "stimuli/img_bwo9g.png"
$ws.Range("M20").Value = 64.81818181818181
$ws.Range("N20").Value = 42.36363636363637
$ws.Range("O20").Value = 53.59090909090909
$ws.Range("P20").Value = 33
$ws.Range("Q20").Value = 4
$ws.Range("R20").Value = 4
$ws.Range("S20").Value = 4
$ws.Range("F21").Value = 87
$ws.Range("L21").Value = "stimuli/img_c79r7.png"
$ws.Range("M21").Value = 56.26470588235294
$ws.Range("N21").Value = 34.26470588235294
$ws.Range("O21").Value = 45.26470588235294
$ws.Range("P21").Value = 34
$ws.Range("Q21").Value = 2
$ws.Range("R21").Value = 2
$ws.Range("S21").Value = 2
$ws.Range("F22").Value = 88
$ws.Range("H22").Value = "bedrooms"
$ws.Range("I22").Value = "distractor"
$ws.Range("K22").Value = "f"
$ws.Range("L22").Value = "stimuli/img_ca8kd.png"
$ws.Range("M22").Value = 92.05405405405405
$ws.Range("N22").Value = 73.02702702702703
$ws.Range("O22").Value = 82.54054054054055
$ws.Range("P22").Value = 37
$ws.Range("Q22").Value = 10
$ws.Range("R22").Value = 10
$ws.Range("S22").Value = 10
$ws.Range("F23").Value = 89
$ws.Range("H23").Value = "kitchens"
$ws.Range("I23").Value = "target"
$ws.Range("K23").Value = "j"
$ws.Range("L23").Value = "stimuli/img_ifebc.png"
$ws.Range("M23").Value = 84
$ws.Range("N23").Value = 65.88235294117646
$ws.Range("O23").Value = 74.94117647058823
$ws.Range("P23").Value = 34
$ws.Range("Q23").Value = 10
$ws.Range("R23").Value = 9
$ws.Range("S23").Value = 9
$ws.Range("F24").Value = 90
$ws.Range("H24").Value = "kitchens"
$ws.Range("I24").Value = "target"
$ws.Range("K24").Value = "j"
$ws.Range("L24").Value = "stimuli/img_p659z.png"
$ws.Range("M24").Value = 84.21621621621621
$ws.Range("N24").Value = 65.37837837837837
$ws.Range("O24").Value = 74.79729729729729
$ws.Range("Q24").Value = 9
$ws.Range("R24").Value = 9
$ws.Range("S24").Value = 9
$ws.Range("F25").Value = 91
$ws.Range("L25").Value = "stimuli/img_ua9bs.png"
$ws.Range("M25").Value = 82
$ws.Range("N25").Value = 62.23333333333333
$ws.Range("O25").Value = 72.11666666666667
$ws.Range("P25").Value = 30
$ws.Range("Q25").Value = 9
$ws.Range("R25").Value = 9
$ws.Range("S25").Value = 9
$ws.Range("F26").Value = 92
$ws.Range("F27").Value = 93
$ws.Range("H27").Value = "living_rooms"
$ws.Range("I27").Value = "distractor"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_xzyzy.png"
$ws.Range("M27").Value = 85.37209302325581
$ws.Range("N27").Value = 68.90697674418605
$ws.Range("O27").Value = 77.13953488372093
$ws.Range("P27").Value = 43
$ws.Range("Q27").Value = 9
$ws.Range("R27").Value = 9
$ws.Range("S27").Value = 9
